$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.297.46'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.707.98'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '587.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.80'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.709.41'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.22%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.514'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.157'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.34'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.444'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.72%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000257'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -7.62%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.35'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.333.09'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.708.48'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.345.16'
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.16'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.93'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.62%  '
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.49'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '460.77'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.690'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.18'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('E25').Value = '  -12.74%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.78'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.96%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.11'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.63%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.857.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.80'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.23'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.67%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.16'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '29.26'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '8.87'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.659.13'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.100'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.35'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -13.05%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.985'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.134'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.31%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.65'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.23%  '
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.300'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.32%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.41'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.88'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.06%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '44.89'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.05%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '143.78'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '386.78'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.43%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0341'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.91%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '24.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.27%  '
